# Weekly data refresh for "Bruselas (repollito)" sheet:
# a new observation row is inserted at sheet row 29 (pushing every
# existing row from 29..113 down by one, to 30..114), and the sheet's
# used-range dimension grows from A1:R113 to A1:R114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 29; everything below shifts down.
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record.
$ws.Cells.Item(29, 1).Value  = 9
$ws.Cells.Item(29, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(29, 3).Value  = "Metropolitana"
$ws.Cells.Item(29, 4).Value  = 45148
$ws.Cells.Item(29, 5).Value  = 13
$ws.Cells.Item(29, 6).Value  = 100112035
$ws.Cells.Item(29, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(29, 8).Value  = "Sin especificar"
$ws.Cells.Item(29, 9).Value  = "Primera"
$ws.Cells.Item(29, 10).Value = 52
$ws.Cells.Item(29, 11).Value = 18000
$ws.Cells.Item(29, 12).Value = 19000
$ws.Cells.Item(29, 13).Value = 18500
$ws.Cells.Item(29, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(29, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(29, 16).Value = 1233
$ws.Cells.Item(29, 17).Value = 15
$ws.Cells.Item(29, 18).Value = "Hortaliza"
